# rerun with small cleans; added var with total persons (tot_pp)
# Table values were recomputed against a slightly different total-persons
# denominator, shifting twenty counts by +/-1, and the whole data block
# (C4:F35) picked up center vertical alignment in addition to its existing
# right alignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Updated counts (rerun with the corrected tot_pp variable) ---
$valueUpdates = @{
    "D4"  = "41659 (47.3)"
    "D5"  = "46326 (52.7)"
    "E5"  = "40308 (52.7)"
    "D6"  = "20017 (22.8)"
    "E6"  = "13509 (17.7)"
    "D7"  = "25100 (28.5)"
    "E11" = "46282 (60.5)"
    "E15" = "62513 (81.7)"
    "D17" = "52824 (60.0)"
    "E17" = "53579 (70.0)"
    "D19" = "7660 (8.7)"
    "D22" = "45383 (51.6)"
    "E22" = "31332 (40.9)"
    "D23" = "23857 (27.1)"
    "D26" = "5232 (5.9)"
    "D31" = "22240 (25.3)"
    "E31" = "20097 (26.3)"
    "D33" = "22695 (25.8)"
    "D34" = "39055 (44.4)"
    "E34" = "49217 (64.3)"
}

foreach ($ref in $valueUpdates.Keys) {
    $ws.Range($ref).Value = $valueUpdates[$ref]
}

# --- 2. Re-format the whole data block with vertical centering ---
$dataRange = $ws.Range("C4:F35")
$dataRange.HorizontalAlignment = -4152   # xlRight (unchanged, kept explicit)
$dataRange.VerticalAlignment = -4108     # xlCenter (new)

# --- 3. Selection left where the author was last working ---
$ws.Range("C15:F16").Select()
